$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 0.083
$ws.Range("M2").Value = 44.1
$ws.Range("N2").Value = 0.9363057324840764
$ws.Range("O2").Value = 531.3253012048193
$ws.Range("S2").Value = 44.1
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 0.222
$ws.Range("V2").Value = 0.004713375796178344
$ws.Range("W2").Value = 0.0166
$ws.Range("X2").Value = 0.04537944222634582
$ws.Range("Y2").Value = -0.02877944222634582
$ws.Range("AA2").Value = -0.2345786272806256
$ws.Range("AB2").Value = 0.04537944222634582
$ws.Range("AC2").Value = -0.2799580695069714
$ws.Range("AG2").Value = -0.222
$ws.Range("AJ2").Value = -0.004735696915397414
$ws.Range("AK2").Value = -0.04646295521138552
$ws.Range("AM2").Value = -1.03
$ws.Range("AQ2").Value = 1.048543689320388
$ws.Range("K3").Value = 0.083
$ws.Range("M3").Value = 44.1
$ws.Range("N3").Value = 0.9363057324840764
$ws.Range("O3").Value = 531.3253012048193
$ws.Range("S3").Value = 44.1
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 0.222
$ws.Range("V3").Value = 0.004713375796178344
$ws.Range("W3").Value = 0.0166
$ws.Range("X3").Value = 0.04537944222634582
$ws.Range("Y3").Value = -0.02877944222634582
$ws.Range("AA3").Value = -0.2345786272806256
$ws.Range("AB3").Value = 0.04537944222634582
$ws.Range("AC3").Value = -0.2799580695069714
$ws.Range("AG3").Value = -0.222
$ws.Range("AJ3").Value = -0.004735696915397414
$ws.Range("AK3").Value = -0.04646295521138552
$ws.Range("AM3").Value = -1.03
$ws.Range("AQ3").Value = 1.048543689320388
